# Flip the sign (negative -> positive) of the values in column G
# for the specified rows on the "Cost" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cost")

$rows = @(45, 78, 86, 140, 141, 148, 273, 274, 293, 294, 295, 296, 297, 298, 299, 300, 301, 302, 303, 304, 305, 306, 307, 308, 309, 310, 311, 312, 313, 314, 315, 316, 317, 318, 319, 320, 321, 322, 323, 324)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)  # Column G is the 7th column
    $cell.Value2 = -1 * $cell.Value2
}
